$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = '1240LT25Q0061'
$ws.Range("B2").Value = 'Kyle Canyon Bunkhouse HVAC Replacement'
$ws.Range("C2").Value = 'Project Name: Kyle Canyon Bunkhouse HVAC Replacement Project Description: This work consists of replacement of an HVAC system at the Kyle Canyon Bunkhouse in the Spring Mountains National Recreation area in the Humboldt-Toiyabe National Forest. The project including removal and disposal of existing furnace, cooling coil and condensing units, in addition to the installation of a new HVAC system with a high altitude conversion kit and liquid propane conversion kit. Project Location: Spring Mountains Recreation Area Bunkhouse (36°16''17.67"N 115°33''46.84"W) Clark County, Nevada Pre-Proposal Site Visit Scheduled? No, see Section L. If interested in a site visit, contact the COR to request. (Dakota.Krist@usda.gov) Schedule Note: The period of performance dates are intended as a timeframe for when the project should be completed (i.e. between October 2025 and March 2026). Once work begins, the awardee will have 30 calendar days to complete the work (Contract time). Submit proposals to Westley.Bisson@usda.gov. Do not mail proposals.'
$ws.Range("D2").Value = 'Oct 10, 2025 05:00 PM MDT'
$ws.Range("E2").Value = 'Sep 12, 2025'
$ws.Range("F2").Value = 'Sep 12, 2025'
$ws.Range("G2").Value = 'Solicitation'
$ws.Range("H2").Value = 'Westley Bisson'
$ws.Range("I2").Value = 'Active'
$ws.Range("J2").Value = 'nan'
$ws.Range("K2").Value = 'nan'
$ws.Range("L2").Value = 'nan'
$ws.Range("M2").Value = 'Plumbing, Heating, and Air-Conditioning Contractors'
$ws.Range("N2").Value = 'MAINT/REPAIR/REBUILD OF EQUIPMENT- REFRIGERATION, AIR CONDITIONING, AND AIR CIRCULATING EQUIPMENT'
$ws.Range("O2").Value = 'nan'
$ws.Range("P2").Value = 'Total Small Business Set-Aside (FAR 19.5)'

# --- Row 3 ---
$ws.Range("A3").Value = '30471PR250000008'
$ws.Range("B3").Value = 'Repair and Install HVAC at USCG STA Morro Bay'
$ws.Range("C3").Value = 'See the attached document for information. Direct all questions to the contracting officer and contract specialist (shawn.t.jenkins@uscg.mil and jesus.c.yutig.civ@uscg.mil). All offers must be emailed to the contracting officer and contract specialist by the closing date and time.'
$ws.Range("D3").Value = 'Sep 19, 2025 03:00 PM EDT'
$ws.Range("E3").Value = 'Sep 13, 2025'
$ws.Range("F3").Value = 'Sep 13, 2025'
$ws.Range("G3").Value = 'Solicitation'
$ws.Range("H3").Value = 'Shawn Jenkins | Jesse Yutig'
$ws.Range("I3").Value = 'Active'
$ws.Range("J3").Value = 'nan'
$ws.Range("K3").Value = 'nan'
$ws.Range("L3").Value = 'nan'
$ws.Range("M3").Value = 'Plumbing, Heating, and Air-Conditioning Contractors'
$ws.Range("N3").Value = 'MAINTENANCE OF OTHER ADMINISTRATIVE FACILITIES AND SERVICE BUILDINGS'
$ws.Range("O3").Value = 'nan'
$ws.Range("P3").Value = 'Total Small Business Set-Aside (FAR 19.5)'

# --- Row 4 ---
$ws.Range("A4").Value = '30890PR250000020'
$ws.Range("B4").Value = 'Replace HVAC at USCG STA Rio Vista'
$ws.Range("C4").Value = 'See the attached document for information. Direct all questions to the contracting officer and contract specialist (shawn.t.jenkins@uscg.mil and jesus.c.yutig.civ@uscg.mil). All offers must be emailed to the contracting officer and contract specialist by the closing date and time.'
$ws.Range("D4").Value = 'Sep 19, 2025 03:00 PM EDT'
$ws.Range("E4").Value = 'Sep 13, 2025'
$ws.Range("F4").Value = 'Sep 13, 2025'
$ws.Range("G4").Value = 'Solicitation'
$ws.Range("H4").Value = 'Shawn Jenkins'
$ws.Range("I4").Value = 'Active'
$ws.Range("J4").Value = 'nan'
$ws.Range("K4").Value = 'nan'
$ws.Range("L4").Value = 'nan'
$ws.Range("M4").Value = 'Plumbing, Heating, and Air-Conditioning Contractors'
$ws.Range("N4").Value = 'MAINTENANCE OF OTHER ADMINISTRATIVE FACILITIES AND SERVICE BUILDINGS'
$ws.Range("O4").Value = 'nan'
$ws.Range("P4").Value = 'Total Small Business Set-Aside (FAR 19.5)'

# --- Row 5 ---
$ws.Range("A5").Value = '70Z028-25-Q-ASMHVACMAINTENANCE'
$ws.Range("B5").Value = 'USCG Air Station Miami FY26 HVAC Maintenance'
$ws.Range("C5").Value = 'This combined synopsis/solicitation for services is prepared in accordance with the format in subpart 12.6 of the Federal Acquisition Regulation (FAR) as supplemented with additional information included in this notice. This synopsis/solicitation is being advertised as Small Business Set-Aside only. This announcement constitutes the only solicitation; no separate written solicitation will be issued. The solicitation number is 70Z028-25-Q-ASMHVACMAINTENANCE Applicable North American Industry Classification Standard (NAICS) codes are: NAICS 238220 - Plumbing, Heating, and Air-Conditioning Contractors Size standard: $19.0 million This requirement is for a firm fixed price services contract. The contract will be awarded using simplified acquisition procedures in accordance with FAR part 12 and FAR Part 13 The Government will award a contract resulting from this solicitation to the responsible offeror whose offer conforming to the solicitation will be most advantageous to the Government, LOWEST PRICE, TECHNICALLY ACCEPTABLE (LPTA). When submitting your proposals, request a price breakdown of the following: a. The cost of monthly HVAC maintenance as outlined in the Statement of Work (SOW), divided into 12 equal monthly payments. b. The annual cost for the Base Year plus 03 Option Years, with pricing for each option year broken down into 12 equal monthly payments. Prospective contractors must also provide, along with your proposal, past performance documents which will be used to determine expertise to perform the work requested. Contract will be awarded to the vendor who can perform the work outlined in the provided PWS in the time frame needed by unit. Award date: 05 business days after close of solicitation, 17 September 2025 Quotes are to be received no later than close of business 11 a.m. on 17 September 2025. Quotes can be email to: Arlett.Maj@uscg.mil Telephone requests of quotes will not be accepted, a formal notice of changes (if applicable) will be issued in Sam.Gov Wage determinations: Service Contract Act WD #: 2015-4543 County: Dade Rev: 32. Date: 07/08/2025 SCOPE OF WORK: To provide Air Station Miami HVAC Maintenane services required on a regular basis for the Air-Cooled Chiller on the NW side of the hangar building, the Water-Cooled Chiller adjacent to the galley building, and the Cooling tower adjacent to the Water-Cooled Chiller at Coast Guard Air Station Miami, Opa Locka, Florida, as defined in the Scope of Work. Location: USCG Air Station Miami 14750 NW 44th Ct Opa Locka, FL 33054 Period of Performance: Performance of work is expected to commence October 1, 2025 Period of Performance: October 1, 2025 – September 30, 2026 52.217-9 Option to Extend the Term of the Contract The Government may extend the term of this contract by written notice to the Contractor provided that the Government gives the Contractor a preliminary written notice of its intent to extend at least __60___ days before the contract expires. The preliminary notice does not commit the Government to an extension. (b) If the Government exercises this option, the extended contract shall be considered to include this option clause. USCG ASM HVAC MAINTENANCE (c) The total duration of this contract, including the exercising of any options under this clause, shall not exceed 48 months. Site visit: It is recommended but not mandatory that a site visit is performed at the Unit to verify the requirements within the PWS and speak with the unit POC for this service. Vendors shall contact: POC: LT Thomas Jones Email: Thomas.M.Jones@uscg.mil Phone: (443) 262-2838 Q&A’s: Questions concerning the work requested must be sent to Arlett.Maj@uscg.mil by COB September 9, 2025. These Q&A’s will be answered and posted to this solicitation as an amendment to solicitation prior to close. 52.212-1 -- Instructions to Offerors -- Commercial Items. (DEVIATION 2018-O0018) Instructions to Offerors -- Commercial Items (June 2020) 52.212-3 Offeror Representations and Certifications -- Commercial Items (June 2020) Work hours: Monday through Friday, 0800 – 1500 (8 am – 3 pm) See attached applicable FAR Clauses by reference. Offerors may obtain a full text version of these clauses electronically at: www.arnet.gov/far The following Homeland Security Acquisition Regulations (HSAR) are incorporated as addenda to this solicitation: HSAR 3052.209-70 Prohibition on Contracts with Corporate Expatriates (Jun 2006); Contracting Officers Technical Representative (Dec 2003) (COTR) HSAR 3052.237-72. Copies of HSAR clauses may be obtained electronically at http://www.dhs.gov. Request Company’s tax ID information and UEI number. Vendors providing an offer must be registered in SAM (System for Award Management) prior to close of this solicitation. https://www.sam.gov/portal/public/SAM/. The vendors SAM registration must be in an “ACTIVE” status prior to award. Vendors must ensure that the above listed NAICS codes are listed on their SAM registration to perform this type of service. The Government reserves the right to cancel this solicitation at any time.'
$ws.Range("D5").Value = 'Sep 17, 2025 11:00 AM EDT'
$ws.Range("E5").Value = 'Sep 4, 2025'
$ws.Range("F5").Value = 'Sep 4, 2025'
$ws.Range("G5").Value = 'Combined Synopsis/Solicitation'
$ws.Range("H5").Value = 'SK1 Arlett Maj | SK1 Alera Nash'
$ws.Range("I5").Value = 'Active'
$ws.Range("J5").Value = 'nan'
$ws.Range("K5").Value = 'nan'
$ws.Range("L5").Value = 'nan'
$ws.Range("M5").Value = 'Plumbing, Heating, and Air-Conditioning Contractors'
$ws.Range("N5").Value = 'MAINT/REPAIR/REBUILD OF EQUIPMENT- REFRIGERATION, AIR CONDITIONING, AND AIR CIRCULATING EQUIPMENT'
$ws.Range("O5").Value = 'nan'
$ws.Range("P5").Value = 'Total Small Business Set-Aside (FAR 19.5)'

# Remove old rows 6 and 7 (shift used range up to A1:P5)
$ws.Range("A6:P7").EntireRow.Delete()

